$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.000.26"
$ws.Range("E2").Value = "'  +0.33%  "
$ws.Range("D3").Value = "'1.642.49"
$ws.Range("E3").Value = "'  +0.03%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D5").Value = "'215.16"
$ws.Range("E5").Value = "'  -0.24%  "
$ws.Range("D6").Value = "'0.5073"
$ws.Range("E6").Value = "'  -0.36%  "
$ws.Range("E7").Value = "'  -0.34%  "
$ws.Range("D8").Value = "'0.2583"
$ws.Range("E8").Value = "'  +0.53%  "
$ws.Range("D9").Value = "'0.06371"
$ws.Range("E9").Value = "'  -0.37%  "
$ws.Range("D10").Value = "'19.88"
$ws.Range("D11").Value = "'0.07730"
$ws.Range("E11").Value = "'  -0.66%  "
$ws.Range("D12").Value = "'4.303"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("D13").Value = "'1.634.94"
$ws.Range("E13").Value = "'  -0.73%  "
$ws.Range("D14").Value = "'0.5490"
$ws.Range("E14").Value = "'  +0.53%  "
$ws.Range("D15").Value = "'0.0₅7760"
$ws.Range("E15").Value = "'  -1.20%  "
$ws.Range("D16").Value = "'64.44"
$ws.Range("E16").Value = "'  -0.32%  "
$ws.Range("D17").Value = "'26.028.04"
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "'  -0.09%  "
$ws.Range("D19").Value = "'198.03"
$ws.Range("E19").Value = "'  +0.19%  "
$ws.Range("D20").Value = "'4.472"
$ws.Range("E20").Value = "'  +0.83%  "
$ws.Range("D21").Value = "'9.992"
$ws.Range("E21").Value = "'  +0.19%  "
$ws.Range("D22").Value = "'6.161"
$ws.Range("E22").Value = "'  +1.90%  "
$ws.Range("E23").Value = "'  -0.53%  "
$ws.Range("D24").Value = "'1.895"
$ws.Range("E24").Value = "'  +0.90%  "
$ws.Range("D25").Value = "'142.65"
$ws.Range("E25").Value = "'  +1.51%  "
$ws.Range("E26").Value = "'  +9.71%  "
$ws.Range("D27").Value = "'6.884"
$ws.Range("E27").Value = "'  -0.23%  "
$ws.Range("D28").Value = "'15.68"
$ws.Range("E28").Value = "'  -0.32%  "
$ws.Range("D29").Value = "'1.243"
$ws.Range("E29").Value = "'  +0.22%  "
$ws.Range("D30").Value = "'0.04903"
$ws.Range("E30").Value = "'  -2.34%  "
$ws.Range("D31").Value = "'3.286"
$ws.Range("E31").Value = "'  +0.67%  "
$ws.Range("D32").Value = "'3.215"
$ws.Range("D33").Value = "'1.559"
$ws.Range("E33").Value = "'  +1.06%  "
$ws.Range("D34").Value = "'2.378"
$ws.Range("E34").Value = "'  +0.52%  "
$ws.Range("D35").Value = "'0.9202"
$ws.Range("E35").Value = "'  +2.88%  "
$ws.Range("D36").Value = "'2.569"
$ws.Range("E36").Value = "'  -0.77%  "
$ws.Range("D37").Value = "'0.5567"
$ws.Range("E37").Value = "'  +1.03%  "
$ws.Range("D38").Value = "'1.112.75"
$ws.Range("E38").Value = "'  -1.77%  "
$ws.Range("D39").Value = "'0.01570"
$ws.Range("E39").Value = "'  +0.76%  "
$ws.Range("E40").Value = "'  -0.51%  "
$ws.Range("D41").Value = "'5.623"
$ws.Range("E41").Value = "'  -0.27%  "
$ws.Range("D42").Value = "'0.8049"
$ws.Range("E42").Value = "'  -1.35%  "
$ws.Range("D43").Value = "'98.77"
$ws.Range("E43").Value = "'  -1.20%  "
$ws.Range("B44").Value = "'RocketPoolETH"
$ws.Range("C44").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "'1.780.96"
$ws.Range("E44").Value = "'  +0.10%  "
$ws.Range("B45").Value = "'BabyDogeCoin"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.0₈119"
$ws.Range("E45").Value = "'  -7.67%  "
$ws.Range("D46").Value = "'0.4532"
$ws.Range("E46").Value = "'  +0.07%  "
$ws.Range("D47").Value = "'55.39"
$ws.Range("E47").Value = "'  +0.82%  "
$ws.Range("E48").Value = "'  +0.16%  "
$ws.Range("D50").Value = "'7.614"
$ws.Range("E50").Value = "'  +2.91%  "
$ws.Range("D51").Value = "'1.005"
$ws.Range("E51").Value = "'  -0.13%  "
